$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new column B for "status_label", shifting existing columns B..I to C..J
$ws.Columns("B").Insert()

# Update the dimension-affecting header
$ws.Cells.Item(1, 2).Value = "status_label"

# Map each row status emoji (col A) to its French color label, written into col B
$statusMap = @{
    "🟥" = "rouge"
    "🟩" = "vert"
    "🟧" = "orange"
}
for ($r = 2; $r -le 14; $r++) {
    $statusValue = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 2).Value = $statusMap[$statusValue]
}

# The source export re-sorted a couple of rows by NCTId; realign the shifted
# NCTId / title / acronym columns (C, F, G) to match for rows 9<->10 and 12<->13
$ws.Cells.Item(9, 3).Value = "NCT05136586"
$ws.Cells.Item(9, 6).Value = "Effects of Two Stress Management Procedures on Performances During Objective Structured Clinical Examination (OSCE) for Medical Students : Relaxing Breathing Combined With Biofeedback or Meditative Stimulation : ECOSTRESS Study"
$ws.Cells.Item(9, 7).Value = "ECOSTRESS"
$ws.Cells.Item(10, 3).Value = "NCT05390879"
$ws.Cells.Item(10, 6).Value = "Influence of Meditation on Stress and Rumination Following Objective Structured Clinical Examination (OSCE)"
$ws.Cells.Item(10, 7).Value = ""
$ws.Cells.Item(12, 3).Value = "NCT05393219"
$ws.Cells.Item(12, 6).Value = "Effects of Preventive Physiological and Psychological Interventions on Performances During Objective Structured Clinical Examination (OSCE) for Medical Students: Cardiac Biofeedback, Mindfulness, or Inner Resources Mobilization"
$ws.Cells.Item(12, 7).Value = ""
$ws.Cells.Item(13, 3).Value = "NCT05619081"
$ws.Cells.Item(13, 6).Value = "Recovery Napping Protocol for Anesthesiologist Performance"
$ws.Cells.Item(13, 7).Value = "R-NAP"
